$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simplify the "Absolute position" values in column B: strip the
# "gi|<number>|ref|" prefix, keeping just "<accession>-<position>".
# (Index 0 of this array lines up with worksheet row 2, index 1 with row 3, etc.)
$newValues = @(
    "NC_017251.1-213522",
    "NC_017250.1-1173757",
    "NC_017250.1-241603",
    "NC_017250.1-1013589",
    "NC_017250.1-1025446",
    "NC_017250.1-1020241",
    "NC_017250.1-1175324",
    "NC_017250.1-291907",
    "NC_017250.1-1197913",
    "NC_017250.1-1016576",
    "NC_017250.1-264518",
    "NC_017251.1-1146371",
    "NC_017250.1-1014292",
    "NC_017250.1-1048661",
    "NC_017250.1-1072841",
    "NC_017250.1-1099151",
    "NC_017251.1-549890",
    "NC_017250.1-1167451",
    "NC_017251.1-1184227",
    "NC_017251.1-1245156",
    "NC_017250.1-1011109",
    "NC_017250.1-1035261",
    "NC_017250.1-1147663",
    "NC_017251.1-994359",
    "NC_017250.1-1195126",
    "NC_017251.1-1088654"
)

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}

# Apply a text number format to column B's data rows (mirrors the new
# cellXfs entry created for these cells).
$rng = $ws.Range("B2:B27")
$rng.NumberFormat = "@"

# Update the active selection to match the author's final cursor position.
$ws.Range("K25").Select()
